# Applies the changes described in the commit:
# 1) openpyexcel was replaced with pandas
# 2) word counters was added

$wb = $excel.ActiveWorkbook
$settings = $wb.Worksheets.Item("settings")
$data = $wb.Worksheets.Item("875771161")

# --- settings sheet ---
# Update the list of "parser" years/tokens in D1
$settings.Range("D1").Value = "2032_2030_2029_2028_2026_2025"

# Column C width (bestFit-like, custom width)
$settings.Columns.Item(3).ColumnWidth = 10

# Update selection on the settings sheet
$settings.Range("K16").Select()

# --- data sheet (875771161) ---
# Remove old rows 10-34 (keep rows 1-9, which are overwritten below)
$data.Rows("10:34").Delete()

# Row 1: pack id + name (drop old D1/E1 "parser:" counter cells)
# A1 keeps the text "4" (leading apostrophe keeps it as text, not a number)
$data.Range("A1").Value = "'4"
$data.Range("B1").Value = "<-pack's name"
$data.Range("D1").Value = ""
$data.Range("E1").Value = ""

# Row 2: headers (unchanged: front / back / description)
$data.Range("A2").Value = "front"
$data.Range("B2").Value = "back"
$data.Range("C2").Value = "description"

# Rows 3-9: front/back numbers (stored as text, like the originals),
# a date (word-counter run date) and a word counter of 1
$dateSerial = 45592
for ($i = 0; $i -lt 7; $i++) {
    $row = 3 + $i
    $data.Cells.Item($row, 1).Value = "'" + [string]($i + 1)
    $data.Cells.Item($row, 2).Value = "'" + [string]($i + 4)
    $data.Cells.Item($row, 3).Value = $dateSerial
    $data.Cells.Item($row, 4).Value = 1
}

# Default row height for the data sheet
$data.Rows.DefaultRowHeight = 15

# Update selection on the data sheet
$data.Range("A1").Select()

# Make the data sheet the active tab
$data.Activate()
